$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.168.43'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '1.849.13'
$ws.Range("E3").Value = '  -0.88%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.20'
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4694'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2890'
$ws.Range("E8").Value = '  +1.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06535'
$ws.Range("E9").Value = '  +0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.73'
$ws.Range("E10").Value = '  +2.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07948'
$ws.Range("E11").Value = '  +1.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.43'
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("D13").Value = '1.852.80'
$ws.Range("E13").Value = '  -0.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.081'
$ws.Range("E14").Value = '  -0.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6738'
$ws.Range("E15").Value = '  +0.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '266.91'
$ws.Range("E16").Value = '  -4.58%  '

$ws.Range("D17").Value = '30.133.32'
$ws.Range("E17").Value = '  -0.43%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.58'
$ws.Range("E18").Value = '  +7.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007591'
$ws.Range("E19").Value = '  +4.57%  '

$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").Value = '2.100.76'
$ws.Range("E21").Value = '  -0.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("E23").Value = '  -5.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.130'
$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '166.63'
$ws.Range("E25").Value = '  +0.86%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.138'
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("E27").Value = '  -1.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.924'
$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.390'
$ws.Range("E29").Value = '  +0.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09824'
$ws.Range("E30").Value = '  +1.96%  '

$ws.Range("E31").Value = '  -0.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.272'
$ws.Range("E32").Value = '  -2.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.995'
$ws.Range("E33").Value = '  -2.37%  '

$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.115'
$ws.Range("E35").Value = '  -0.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6967'
$ws.Range("E36").Value = '  -1.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.707'
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01862'
$ws.Range("E38").Value = '  +0.52%  '

$ws.Range("E39").Value = '  +2.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.328'
$ws.Range("E40").Value = '  +1.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.01'
$ws.Range("E41").Value = '  -0.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.928'
$ws.Range("E42").Value = '  -0.74%  '

$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8354'
$ws.Range("E44").Value = '  -1.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.02'
$ws.Range("E45").Value = '  -0.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4118'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '938.12'
$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.115'
$ws.Range("E48").Value = '  -0.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.987'
$ws.Range("E49").Value = '  -2.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.78'
$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05653'
$ws.Range("E51").Value = '  +0.45%  '
Write-Host "Updated cryptos list"
